# Updates the cryptos list: refresh Price (column D) and Volume(1h) (column E)
# values for rows 2-51 on the active worksheet, matching the latest scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "37.763.63"; DForceText = $false; E = "-1.25%" },
    @{ Row = 3; D = "2.031.79"; DForceText = $false; E = "-1.69%" },
    @{ Row = 4; D = $null; DForceText = $false; E = "-0.10%" },
    @{ Row = 5; D = "227.50"; DForceText = $true; E = "-1.26%" },
    @{ Row = 6; D = $null; DForceText = $false; E = "-1.72%" },
    @{ Row = 7; D = "59.89"; DForceText = $true; E = "-2.22%" },
    @{ Row = 8; D = $null; DForceText = $false; E = "+0.13%" },
    @{ Row = 9; D = $null; DForceText = $false; E = "-2.97%" },
    @{ Row = 10; D = $null; DForceText = $false; E = "+2.32%" },
    @{ Row = 11; D = $null; DForceText = $false; E = "-0.19%" },
    @{ Row = 12; D = "14.66"; DForceText = $true; E = "-1.49%" },
    @{ Row = 13; D = "2.331.67"; DForceText = $false; E = "-1.83%" },
    @{ Row = 14; D = "21.04"; DForceText = $true; E = "-1.18%" },
    @{ Row = 15; D = "0.766"; DForceText = $true; E = "+0.10%" },
    @{ Row = 16; D = $null; DForceText = $false; E = "-2.53%" },
    @{ Row = 17; D = "2.025.73"; DForceText = $false; E = "-4.76%" },
    @{ Row = 18; D = "37.713.18"; DForceText = $false; E = "-1.24%" },
    @{ Row = 19; D = "69.51"; DForceText = $true; E = "-0.93%" },
    @{ Row = 20; D = "5.89"; DForceText = $true; E = "-6.08%" },
    @{ Row = 21; D = $null; DForceText = $false; E = "-1.53%" },
    @{ Row = 22; D = "223.79"; DForceText = $true; E = "-1.01%" },
    @{ Row = 23; D = $null; DForceText = $false; E = "+0.13%" },
    @{ Row = 24; D = "2.37"; DForceText = $true; E = "-2.08%" },
    @{ Row = 25; D = "2.27"; DForceText = $true; E = "+1.50%" },
    @{ Row = 26; D = "168.27"; DForceText = $true; E = "+1.18%" },
    @{ Row = 27; D = $null; DForceText = $false; E = "+1.07%" },
    @{ Row = 28; D = $null; DForceText = $false; E = "-4.16%" },
    @{ Row = 29; D = "18.76"; DForceText = $true; E = "-1.06%" },
    @{ Row = 30; D = $null; DForceText = $false; E = "-4.23%" },
    @{ Row = 31; D = $null; DForceText = $false; E = "+0.89%" },
    @{ Row = 32; D = "2.23"; DForceText = $true; E = "+9.15%" },
    @{ Row = 33; D = $null; DForceText = $false; E = "-4.16%" },
    @{ Row = 34; D = "0.0604"; DForceText = $true; E = "-0.15%" },
    @{ Row = 35; D = $null; DForceText = $false; E = "-2.99%" },
    @{ Row = 36; D = "6.45"; DForceText = $true; E = "+3.02%" },
    @{ Row = 37; D = $null; DForceText = $false; E = "-1.31%" },
    @{ Row = 38; D = "3.41"; DForceText = $true; E = "+2.99%" },
    @{ Row = 39; D = $null; DForceText = $false; E = "+0.05%" },
    @{ Row = 40; D = "18.09"; DForceText = $true; E = "+5.71%" },
    @{ Row = 41; D = "1.537.57"; DForceText = $false; E = "+0.98%" },
    @{ Row = 42; D = $null; DForceText = $false; E = "-0.97%" },
    @{ Row = 43; D = "95.81"; DForceText = $true; E = "-2.49%" },
    @{ Row = 44; D = "2.80"; DForceText = $true; E = "-2.89%" },
    @{ Row = 45; D = "0.0905"; DForceText = $true; E = "-2.47%" },
    @{ Row = 46; D = $null; DForceText = $false; E = "+0.93%" },
    @{ Row = 47; D = $null; DForceText = $false; E = "-3.02%" },
    @{ Row = 48; D = $null; DForceText = $false; E = "-1.75%" },
    @{ Row = 49; D = $null; DForceText = $false; E = "+0.15%" },
    @{ Row = 50; D = "7.10"; DForceText = $true; E = "-0.29%" },
    @{ Row = 51; D = "2.220.96"; DForceText = $false; E = "-1.81%" }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($r, 4)
        if ($u.DForceText) {
            # Values that would otherwise be auto-parsed as numbers (and lose
            # their trailing zeros / formatting) must be forced to text.
            $cell.NumberFormat = "@"
        }
        $cell.Value = $u.D
    }
    $ws.Cells.Item($r, 5).Value = "  " + $u.E + "  "
}
